# Parse countries and country codes from original address
# -> New lab member "Gurpreet Singh Khalsa" (Jammu, India) is appended
#    to the "2021" members table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2021")

# Grow the existing structured table by one row so the table ref /
# autoFilter / dimension all stay consistent, then fill in the new data.
$lo = $ws.ListObjects.Item(1)
[void]$lo.ListRows.Add()

$ws.Range("A11").Value = "Gurpreet Singh Khalsa"
$ws.Range("B11").Value = "Jammu, India"

# Match the selection left behind after typing the new address.
[void]$ws.Range("B11").Select()
